$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 10.12973313815042
$ws.Cells.Item(2, 3).Value = 5.394674703986428
$ws.Cells.Item(2, 5).Value = 16.39776727728987
$ws.Cells.Item(2, 6).Value = 16.86991607391245
$ws.Cells.Item(2, 7).Value = 36.17661620787334
$ws.Cells.Item(2, 8).Value = 15.85790271780439
$ws.Cells.Item(2, 11).Value = 9.372913197012805

$ws.Cells.Item(3, 2).Value = 9.840209347290049
$ws.Cells.Item(3, 3).Value = 5.185536443411548
$ws.Cells.Item(3, 5).Value = 15.47323879823657
$ws.Cells.Item(3, 6).Value = 15.89584955866815
$ws.Cells.Item(3, 7).Value = 35.98560517980675
$ws.Cells.Item(3, 8).Value = 15.88507003811607
$ws.Cells.Item(3, 11).Value = 9.177164971097383

$ws.Cells.Item(4, 2).Value = 9.660893138802878
$ws.Cells.Item(4, 3).Value = 5.051527209673093
$ws.Cells.Item(4, 5).Value = 14.88159279326575
$ws.Cells.Item(4, 6).Value = 15.26997757108489
$ws.Cells.Item(4, 7).Value = 35.88060449216924
$ws.Cells.Item(4, 8).Value = 15.90499443110911
$ws.Cells.Item(4, 11).Value = 9.057374919679134

$ws.Cells.Item(5, 2).Value = 9.587556444854613
$ws.Cells.Item(5, 3).Value = 4.995551070733725
$ws.Cells.Item(5, 5).Value = 14.63473444617097
$ws.Cells.Item(5, 6).Value = 15.008197319934
$ws.Cells.Item(5, 7).Value = 35.84092265894309
$ws.Cells.Item(5, 8).Value = 15.91392575714618
$ws.Cells.Item(5, 11).Value = 9.008734779962442

$ws.Cells.Item(6, 2).Value = 9.575366794977718
$ws.Cells.Item(6, 3).Value = 4.986175070586703
$ws.Cells.Item(6, 5).Value = 14.59340505992377
$ws.Cells.Item(6, 6).Value = 14.96433081551589
$ws.Cells.Item(6, 7).Value = 35.83452159374572
$ws.Cells.Item(6, 8).Value = 15.91545772830635
$ws.Cells.Item(6, 11).Value = 9.000670915526559

$ws.Cells.Item(7, 2).Value = 9.659904992587339
$ws.Cells.Item(7, 3).Value = 5.050777769288202
$ws.Cells.Item(7, 5).Value = 14.87828648917138
$ws.Cells.Item(7, 6).Value = 15.26647399323133
$ws.Cells.Item(7, 7).Value = 35.88005672969609
$ws.Cells.Item(7, 8).Value = 15.90511159947421
$ws.Cells.Item(7, 11).Value = 9.056718132466708

$ws.Cells.Item(8, 2).Value = 10.0302995017162
$ws.Cells.Item(8, 3).Value = 5.323748639585678
$ws.Cells.Item(8, 5).Value = 16.08411518587273
$ws.Cells.Item(8, 6).Value = 16.5399640634477
$ws.Cells.Item(8, 7).Value = 36.10821969813101
$ws.Cells.Item(8, 8).Value = 15.86659467079273
$ws.Cells.Item(8, 11).Value = 9.305380488735675

$ws.Cells.Item(9, 2).Value = 10.73922638865524
$ws.Cells.Item(9, 3).Value = 5.813089659827588
$ws.Cells.Item(9, 5).Value = 18.30270266975359
$ws.Cells.Item(9, 6).Value = 19.00274580682531
$ws.Cells.Item(9, 7).Value = 36.65199336892591
$ws.Cells.Item(9, 8).Value = 15.81695560970081
$ws.Cells.Item(9, 11).Value = 9.792988697588388

$ws.Cells.Item(10, 2).Value = 11.24313493738778
$ws.Cells.Item(10, 3).Value = 6.14287222061152
$ws.Cells.Item(10, 5).Value = 19.93413146704076
$ws.Cells.Item(10, 6).Value = 20.67494806633232
$ws.Cells.Item(10, 7).Value = 37.10838989838404
$ws.Cells.Item(10, 8).Value = 15.79647483276256
$ws.Cells.Item(10, 11).Value = 10.14711646544621

$ws.Cells.Item(11, 2).Value = 11.46752717483802
$ws.Cells.Item(11, 3).Value = 6.286162852228859
$ws.Cells.Item(11, 5).Value = 20.6346565638774
$ws.Cells.Item(11, 6).Value = 21.3917225636224
$ws.Cells.Item(11, 7).Value = 37.32786398841801
$ws.Cells.Item(11, 8).Value = 15.79066916521361
$ws.Cells.Item(11, 11).Value = 10.30651535319929

$ws.Cells.Item(12, 2).Value = 11.55171687547554
$ws.Cells.Item(12, 3).Value = 6.339436956441398
$ws.Cells.Item(12, 5).Value = 20.893979656225
$ws.Cells.Item(12, 6).Value = 21.65686569030329
$ws.Cells.Item(12, 7).Value = 37.4126282800045
$ws.Cells.Item(12, 8).Value = 15.78897847257436
$ws.Cells.Item(12, 11).Value = 10.36656920808739

$ws.Cells.Item(13, 2).Value = 11.53362122875186
$ws.Cells.Item(13, 3).Value = 6.328007594235195
$ws.Cells.Item(13, 5).Value = 20.83839389627077
$ws.Cells.Item(13, 6).Value = 21.60004134736742
$ws.Cells.Item(13, 7).Value = 37.39430005584335
$ws.Cells.Item(13, 8).Value = 15.78931996594191
$ws.Cells.Item(13, 11).Value = 10.35365015652611

$ws.Cells.Item(14, 2).Value = 11.47446966910562
$ws.Cells.Item(14, 3).Value = 6.290565624308556
$ws.Cells.Item(14, 5).Value = 20.65611029275584
$ws.Cells.Item(14, 6).Value = 21.4136618050453
$ws.Cells.Item(14, 7).Value = 37.33480475100635
$ws.Cells.Item(14, 8).Value = 15.79051987837338
$ws.Cells.Item(14, 11).Value = 10.31146252173882

$ws.Cells.Item(15, 2).Value = 11.4381332497636
$ws.Cells.Item(15, 3).Value = 6.267502287049305
$ws.Cells.Item(15, 5).Value = 20.54368239999211
$ws.Cells.Item(15, 6).Value = 21.29868154950795
$ws.Cells.Item(15, 7).Value = 37.29857600529884
$ws.Cells.Item(15, 8).Value = 15.79132107010239
$ws.Cells.Item(15, 11).Value = 10.28557958183019

$ws.Cells.Item(16, 2).Value = 11.2283653886003
$ws.Cells.Item(16, 3).Value = 6.133370710886362
$ws.Cells.Item(16, 5).Value = 19.88751594411908
$ws.Cells.Item(16, 6).Value = 20.62722412089977
$ws.Cells.Item(16, 7).Value = 37.09428110059179
$ws.Cells.Item(16, 8).Value = 15.79692515225143
$ws.Cells.Item(16, 11).Value = 10.13665976007965

$ws.Cells.Item(17, 2).Value = 11.09837582593741
$ws.Cells.Item(17, 3).Value = 6.049347089031314
$ws.Cells.Item(17, 5).Value = 19.47434369432803
$ws.Cells.Item(17, 6).Value = 20.20408069617459
$ws.Cells.Item(17, 7).Value = 36.9719549149282
$ws.Cells.Item(17, 8).Value = 15.80126449199447
$ws.Cells.Item(17, 11).Value = 10.04482177407302

$ws.Cells.Item(18, 2).Value = 11.02315853131434
$ws.Cells.Item(18, 3).Value = 6.0003864874778
$ws.Cells.Item(18, 5).Value = 19.23277842688593
$ws.Cells.Item(18, 6).Value = 19.95656407809808
$ws.Cells.Item(18, 7).Value = 36.90271550172602
$ws.Cells.Item(18, 8).Value = 15.80409069238675
$ws.Cells.Item(18, 11).Value = 9.991842934089611

$ws.Cells.Item(19, 2).Value = 10.99761662178247
$ws.Cells.Item(19, 3).Value = 5.983701288131293
$ws.Cells.Item(19, 5).Value = 19.15031434741686
$ws.Cells.Item(19, 6).Value = 19.87204792380562
$ws.Cells.Item(19, 7).Value = 36.87946595543535
$ws.Cells.Item(19, 8).Value = 15.80510423734048
$ws.Cells.Item(19, 11).Value = 9.973880395840954

$ws.Cells.Item(20, 2).Value = 11.1122607431609
$ws.Cells.Item(20, 3).Value = 6.058357151569787
$ws.Cells.Item(20, 5).Value = 19.51873199044467
$ws.Cells.Item(20, 6).Value = 20.24955283636157
$ws.Cells.Item(20, 7).Value = 36.98486124990664
$ws.Cells.Item(20, 8).Value = 15.80076835278094
$ws.Cells.Item(20, 11).Value = 10.05461470644889

$ws.Cells.Item(21, 2).Value = 11.49186581269206
$ws.Cells.Item(21, 3).Value = 6.301590162451205
$ws.Cells.Item(21, 5).Value = 20.70981262941732
$ws.Cells.Item(21, 6).Value = 21.46857628470577
$ws.Cells.Item(21, 7).Value = 37.3522354954206
$ws.Cells.Item(21, 8).Value = 15.79015363191715
$ws.Cells.Item(21, 11).Value = 10.32386287013237

$ws.Cells.Item(22, 2).Value = 11.7353625022803
$ws.Cells.Item(22, 3).Value = 6.454797470489389
$ws.Cells.Item(22, 5).Value = 21.45358406938259
$ws.Cells.Item(22, 6).Value = 22.22866616901552
$ws.Cells.Item(22, 7).Value = 37.60195134661809
$ws.Cells.Item(22, 8).Value = 15.78617724765146
$ws.Cells.Item(22, 11).Value = 10.49801589771837

$ws.Cells.Item(23, 2).Value = 11.6058510029802
$ws.Cells.Item(23, 3).Value = 6.373560435538284
$ws.Cells.Item(23, 5).Value = 21.05978013657664
$ws.Cells.Item(23, 6).Value = 21.82633154458858
$ws.Cells.Item(23, 7).Value = 37.46781134606844
$ws.Cells.Item(23, 8).Value = 15.7880276977632
$ws.Cells.Item(23, 11).Value = 10.40525319602538

$ws.Cells.Item(24, 2).Value = 11.10598487467604
$ws.Cells.Item(24, 3).Value = 6.054285738859056
$ws.Cells.Item(24, 5).Value = 19.49867658810662
$ws.Cells.Item(24, 6).Value = 20.22900810905287
$ws.Cells.Item(24, 7).Value = 36.97902290442926
$ws.Cells.Item(24, 8).Value = 15.80099162504433
$ws.Cells.Item(24, 11).Value = 10.05018787738265

$ws.Cells.Item(25, 2).Value = 10.54999226189833
$ws.Cells.Item(25, 3).Value = 5.685824452426333
$ws.Cells.Item(25, 5).Value = 17.68583023637704
$ws.Cells.Item(25, 6).Value = 18.34778573295695
$ws.Cells.Item(25, 7).Value = 36.49474301863735
$ws.Cells.Item(25, 8).Value = 15.82759042982474
$ws.Cells.Item(25, 11).Value = 9.661514900567237
